$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# The old row10 "Contact"/"No display for ContactDetail" becomes
# "Jurisdiction"/"United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The old row11 was a duplicate "Contact"/"No display for ContactDetail" row;
# remove it entirely so everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" (sheet2) ---
$elements = $wb.Worksheets.Item("Elements")

# Short / Definition for the root Extension row now show the profile-specific
# text instead of the generic Extension placeholder text.
$elements.Range("K2").Value = "Plan Employee Code"
$elements.Range("L2").Value = "Code indicating whether the member is an employee of the health plan or the dependent of an employee of the health plan"
